$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.232253
$ws.Range("H2").Value = 0.6967589999999999
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 47.890007
$ws.Range("N2").Value = 95.780014
$ws.Range("O2").Value = 0.2671069186198896
$ws.Range("P2").Value = 0.2091536763465644
$ws.Range("Q2").Value = 11.122597795771
$ws.Range("R2").Value = 66.73558677462599
$ws.Range("S2").Value = 0.2671069186198896
$ws.Range("T2").Value = 0.2091536763465644

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.232253
$ws.Range("H3").Value = 0.6967589999999999
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.489942
$ws.Range("N3").Value = 73.46982600000001
$ws.Range("O3").Value = 0.136592858397365
$ws.Range("P3").Value = 0.1604351844054064
$ws.Range("Q3").Value = 5.687862499326
$ws.Range("R3").Value = 51.190762493934
$ws.Range("S3").Value = 0.136592858397365
$ws.Range("T3").Value = 0.1604351844054064

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.232253
$ws.Range("H4").Value = 0.6967589999999999
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.62554866666667
$ws.Range("N4").Value = 61.876646
$ws.Range("O4").Value = 0.1150391719340928
$ws.Range("P4").Value = 0.1351192952518773
$ws.Range("Q4").Value = 4.790345554479333
$ws.Range("R4").Value = 43.11310999031399
$ws.Range("S4").Value = 0.1150391719340928
$ws.Range("T4").Value = 0.1351192952518773

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.232253
$ws.Range("H5").Value = 0.6967589999999999
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.83893633333333
$ws.Range("N5").Value = 65.516809
$ws.Range("O5").Value = 0.1218068518956912
$ws.Range("P5").Value = 0.1430682758601985
$ws.Range("Q5").Value = 5.072158480225665
$ws.Range("R5").Value = 45.64942632203099
$ws.Range("S5").Value = 0.1218068518956912
$ws.Range("T5").Value = 0.1430682758601985

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.232253
$ws.Range("H6").Value = 0.6967589999999999
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 32.40337933333333
$ws.Range("N6").Value = 97.210138
$ws.Range("O6").Value = 0.1807301219771816
$ws.Range("P6").Value = 0.2122766211002732
$ws.Range("Q6").Value = 7.525782060304666
$ws.Range("R6").Value = 67.73203854274199
$ws.Range("S6").Value = 0.1807301219771816
$ws.Range("T6").Value = 0.2122766211002732

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.232253
$ws.Range("H7").Value = 0.6967589999999999
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 32.0437125
$ws.Range("N7").Value = 64.087425
$ws.Range("O7").Value = 0.1787240771757799
$ws.Range("P7").Value = 0.1399469470356803
$ws.Range("Q7").Value = 7.442248359262498
$ws.Range("R7").Value = 44.65349015557499
$ws.Range("S7").Value = 0.1787240771757799
$ws.Range("T7").Value = 0.1399469470356803

Write-Output "Applied Gdf6-Bmpr2 updates (Natmi following Dr Hou advice)"